$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.989.61"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.09%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.896.27"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.75%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9984"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.18%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.8449"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +6.76%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.28"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.05%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9983"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.22%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3305"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +4.67%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "26.75"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.09%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07071"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.78%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08064"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.34%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7605"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.49%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.897.35"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.37%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.268"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.88%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.43"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.65%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.988.94"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.06%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.13"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.92%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.899"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.03%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "245.41"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.01%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007781"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.07%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9977"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.26%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.147.33"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.33%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9978"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.21%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.005"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.38%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1734"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +28.57%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.260"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.63%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.06"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.59%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.92"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.65%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.114"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.82%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.360"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.09%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.516"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.55%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.05890"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +9.74%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.311"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.87%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.091"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.72%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.278"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.13%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7346"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.53%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.709"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.79%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01922"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.18%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.777"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.37%  "

# Row 40
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.12%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "72.67"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.34%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.895"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -4.61%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8471"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.86%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9977"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.36%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.890"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.76%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.68"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.02%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.611"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.42%  "

# Row 48
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.011.86"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +5.45%  "

# Row 49
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.808"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.69%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.047.13"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.19%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "35.99"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.03%  "
